$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Update the celular number and document-user fields in row 2
$ws.Range("B2").Value = 93221452
$ws.Range("D2").Value = "autotest27"

# Left-align the celular column value (new style entry)
$ws.Range("B2").HorizontalAlignment = -4131  # xlLeft

# Reset the view: scroll back to show column A and select A2
$ws.Activate()
$ws.Range("A2").Select()
